$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
  # row 5
$ws.Range("H5").Value = 294
$ws.Range("I5").Value = 367
$ws.Range("J5").Value = 75
$ws.Range("K5").Value = 367
$ws.Range("L5").Value = 75
$ws.Range("M5").Value = -252
$ws.Range("N5").Value = -305
  # row 10
$ws.Range("H10").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("N10").ClearContents()
  # row 11
$ws.Range("H11").Value = 713.4
$ws.Range("I11").Value = 713.4
$ws.Range("K11").Value = 713.4
$ws.Range("M11").Value = -573.4
  # row 33
$ws.Range("H33").Value = 521.0909
$ws.Range("I33").Value = 294.64706
$ws.Range("K33").Value = 294.64706
$ws.Range("M33").Value = -65.64706000000001
  # row 41
$ws.Range("H41").Value = 5402.5
$ws.Range("I41").Value = 5516.375
$ws.Range("K41").Value = 5516.375
$ws.Range("M41").Value = -5076.375
  # row 86
$ws.Range("H86").Value = 214292060
$ws.Range("I86").Value = 200007260
$ws.Range("J86").Value = 250004080
$ws.Range("K86").Value = 200007260
$ws.Range("L86").Value = 250004080
$ws.Range("M86").Value = -200006137
$ws.Range("N86").Value = -250006326
  # row 89
$ws.Range("H89").Value = 214292060
$ws.Range("I89").Value = 200007260
$ws.Range("J89").Value = 250004080
$ws.Range("K89").Value = 1000036300
$ws.Range("L89").Value = 1250020400
$ws.Range("M89").Value = -1000030684
$ws.Range("N89").Value = -1250031632
  # row 100
$ws.Range("H100").Value = 2456.5715
$ws.Range("I100").Value = 1139.4
$ws.Range("K100").Value = 1139.4
$ws.Range("M100").Value = -598.4000000000001
  # row 113
$ws.Range("H113").Value = 3792.75
$ws.Range("I113").Value = 3886.5557
$ws.Range("K113").Value = 3886.5557
$ws.Range("M113").Value = -632.5556999999999
  # row 116
$ws.Range("H116").Value = 33799
$ws.Range("I116").Value = 41248.75
$ws.Range("K116").Value = 41248.75
$ws.Range("M116").Value = -37806.75
  # row 135
$ws.Range("H135").Value = 1140.25
$ws.Range("I135").Value = 1198.3636
$ws.Range("K135").Value = 10785.2724
$ws.Range("M135").Value = -8250.2724
  # row 137
$ws.Range("H137").Value = 3476186.8
$ws.Range("I137").Value = 2818.8333
$ws.Range("K137").Value = 8456.499899999999
$ws.Range("M137").Value = -5906.499899999999
  # row 138
$ws.Range("H138").Value = 2843.29
$ws.Range("I138").Value = 1794.6154
$ws.Range("J138").Value = 2999.9885
$ws.Range("K138").Value = 5383.8462
$ws.Range("L138").Value = 8999.9655
$ws.Range("M138").Value = -243.8462
$ws.Range("N138").Value = -19279.9655
  # row 141
$ws.Range("H141").Value = 4929.6665
$ws.Range("I141").Value = 5043.1763
$ws.Range("K141").Value = 15129.5289
$ws.Range("M141").Value = -9949.528900000001

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
  # row 32
$ws.Range("H32").Value = 18546208
$ws.Range("I32").Value = 18595994
$ws.Range("K32").Value = 18595994
$ws.Range("M32").Value = -18595707
  # row 63
$ws.Range("H63").Value = 3512.5386
$ws.Range("I63").Value = 2044
$ws.Range("K63").Value = 2044
$ws.Range("M63").Value = -1358
  # row 66
$ws.Range("H66").Value = 3512.5386
$ws.Range("I66").Value = 2044
$ws.Range("K66").Value = 10220
$ws.Range("M66").Value = -6788
  # row 74
$ws.Range("H74").Value = 1810.4872
$ws.Range("I74").Value = 1778.4324
$ws.Range("K74").Value = 1778.4324
$ws.Range("M74").Value = -904.4323999999999
  # row 77
$ws.Range("H77").Value = 1810.4872
$ws.Range("I77").Value = 1778.4324
$ws.Range("K77").Value = 8892.162
$ws.Range("M77").Value = -4524.162
  # row 80
$ws.Range("H80").Value = 99999
$ws.Range("I80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("M80").ClearContents()
  # row 83
$ws.Range("H83").Value = 99999
$ws.Range("I83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("M83").ClearContents()
  # row 88
$ws.Range("H88").Value = 1013
$ws.Range("J88").Value = 515.8333
$ws.Range("L88").Value = 515.8333
$ws.Range("N88").Value = -1327.8333
  # row 91
$ws.Range("H91").Value = 1013
$ws.Range("J91").Value = 515.8333
$ws.Range("L91").Value = 515.8333
$ws.Range("N91").Value = -3323.8333
  # row 105
$ws.Range("H105").Value = 112999
$ws.Range("J105").Value = 112999
$ws.Range("L105").Value = 112999
$ws.Range("N105").Value = -119987
  # row 132
$ws.Range("H132").Value = 2164.4375
$ws.Range("I132").Value = 2268.8215
$ws.Range("K132").Value = 6806.4645
$ws.Range("M132").Value = -4276.4645
  # row 134
$ws.Range("H134").Value = 120429
$ws.Range("J134").Value = 120429
$ws.Range("L134").Value = 120429
$ws.Range("N134").Value = -130569

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
  # row 50
$ws.Range("H50").Value = 51978.5
$ws.Range("J50").Value = 51978.5
$ws.Range("L50").Value = 51978.5
$ws.Range("N50").Value = -53126.5
  # row 63
$ws.Range("H63").Value = 70210.55499999999
$ws.Range("J63").Value = 70210.55499999999
$ws.Range("L63").Value = 70210.55499999999
$ws.Range("N63").Value = -71582.55499999999
  # row 66
$ws.Range("H66").Value = 70210.55499999999
$ws.Range("J66").Value = 70210.55499999999
$ws.Range("L66").Value = 210631.665
$ws.Range("N66").Value = -217495.665
  # row 82
$ws.Range("H82").Value = 37314.777
$ws.Range("I82").Value = 27262.143
$ws.Range("J82").Value = 72499
$ws.Range("K82").Value = 27262.143
$ws.Range("L82").Value = 72499
$ws.Range("M82").Value = -26879.143
$ws.Range("N82").Value = -73265
  # row 85
$ws.Range("H85").Value = 37314.777
$ws.Range("I85").Value = 27262.143
$ws.Range("J85").Value = 72499
$ws.Range("K85").Value = 27262.143
$ws.Range("L85").Value = 72499
$ws.Range("M85").Value = -25936.143
$ws.Range("N85").Value = -75151
  # row 92
$ws.Range("H92").Value = 55981.76
$ws.Range("J92").Value = 55981.76
$ws.Range("L92").Value = 55981.76
$ws.Range("N92").Value = -60973.76
  # row 107
$ws.Range("H107").Value = 2382.9048
$ws.Range("I107").Value = 3031.5
$ws.Range("K107").Value = 3031.5
$ws.Range("M107").Value = -1111.5
  # row 134
$ws.Range("H134").Value = 2085791.8
$ws.Range("I134").Value = 2566497.8
$ws.Range("J134").Value = 2732.6667
$ws.Range("K134").Value = 7699493.399999999
$ws.Range("L134").Value = 8198.000100000001
$ws.Range("M134").Value = -7696958.399999999
$ws.Range("N134").Value = -13268.0001

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
  # row 16
$ws.Range("H16").Value = 1644.4062
$ws.Range("I16").Value = 1627.0714
$ws.Range("J16").Value = 1765.75
$ws.Range("K16").Value = 1627.0714
$ws.Range("L16").Value = 1765.75
$ws.Range("M16").Value = -1340.0714
$ws.Range("N16").Value = -2339.75
  # row 22
$ws.Range("H22").Value = 1665.909
$ws.Range("I22").Value = 1980.6666
$ws.Range("K22").Value = 1980.6666
$ws.Range("M22").Value = -1630.6666
  # row 58
$ws.Range("H58").Value = 3097.7546
$ws.Range("I58").Value = 2612.2896
$ws.Range("J58").Value = 4327.6
$ws.Range("K58").Value = 2612.2896
$ws.Range("L58").Value = 4327.6
$ws.Range("M58").Value = -2409.2896
$ws.Range("N58").Value = -4733.6
  # row 105
$ws.Range("H105").Value = 2331.1904
$ws.Range("I105").Value = 2272.75
$ws.Range("K105").Value = 2272.75
$ws.Range("M105").Value = -525.75
  # row 113
$ws.Range("H113").Value = 1644.4062
$ws.Range("I113").Value = 1627.0714
$ws.Range("J113").Value = 1765.75
$ws.Range("K113").Value = 1627.0714
$ws.Range("L113").Value = 1765.75
$ws.Range("M113").Value = 542.9286
$ws.Range("N113").Value = -6105.75
  # row 134
$ws.Range("H134").Value = 2519.276
$ws.Range("I134").Value = 2449.423
$ws.Range("K134").Value = 7348.268999999999
$ws.Range("M134").Value = -4813.268999999999
  # row 136
$ws.Range("H136").Value = 3097.7546
$ws.Range("I136").Value = 2612.2896
$ws.Range("J136").Value = 4327.6
$ws.Range("K136").Value = 7836.8688
$ws.Range("L136").Value = 12982.8
$ws.Range("M136").Value = -5286.8688
$ws.Range("N136").Value = -18082.8
  # row 141
$ws.Range("H141").Value = 653991.2
$ws.Range("J141").Value = 653991.2
$ws.Range("L141").Value = 653991.2
$ws.Range("N141").Value = -664351.2

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
  # row 2
$ws.Range("H2").Value = 335.5
$ws.Range("J2").Value = 428
$ws.Range("L2").Value = 2568
$ws.Range("N2").Value = -2794
  # row 11
$ws.Range("H11").Value = 1666900
$ws.Range("J11").Value = 500
$ws.Range("L11").Value = 1500
$ws.Range("N11").Value = -1780
  # row 22
$ws.Range("H22").Value = 600
$ws.Range("I22").Value = 600
$ws.Range("K22").Value = 1800
$ws.Range("M22").Value = -1631
  # row 26
$ws.Range("H26").Value = 133.5
$ws.Range("I26").Value = 60.2
$ws.Range("J26").Value = 500
$ws.Range("K26").Value = 180.6
$ws.Range("L26").Value = 1500
$ws.Range("M26").Value = 107.4
$ws.Range("N26").Value = -2076
  # row 27
$ws.Range("H27").Value = 600
$ws.Range("I27").Value = 600
$ws.Range("K27").Value = 1800
$ws.Range("M27").Value = -1698
  # row 92
$ws.Range("H92").Value = 1343.7142
$ws.Range("J92").Value = 1965.3334
$ws.Range("L92").Value = 5896.0002
$ws.Range("N92").Value = -8392.0002
  # row 128
$ws.Range("H128").Value = 399998.66
$ws.Range("I128").Value = 399998.66
$ws.Range("K128").Value = 1199995.98
$ws.Range("M128").Value = -1195015.98
  # row 131
$ws.Range("H131").Value = 1675.5358
$ws.Range("J131").Value = 1843.9048
$ws.Range("L131").Value = 5531.7144
$ws.Range("N131").Value = -15611.7144

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
  # row 29
$ws.Range("H29").Value = 2272.7273
$ws.Range("J29").Value = 2272.7273
$ws.Range("L29").Value = 2272.7273
$ws.Range("N29").Value = -2852.7273
  # row 80
$ws.Range("H80").Value = 2128.0833
$ws.Range("I80").Value = 1852
$ws.Range("J80").Value = 2325.2856
$ws.Range("K80").Value = 1852
$ws.Range("L80").Value = 2325.2856
$ws.Range("M80").Value = -854
$ws.Range("N80").Value = -4321.2856
  # row 83
$ws.Range("H83").Value = 2128.0833
$ws.Range("I83").Value = 1852
$ws.Range("J83").Value = 2325.2856
$ws.Range("K83").Value = 9260
$ws.Range("L83").Value = 11626.428
$ws.Range("M83").Value = -4268
$ws.Range("N83").Value = -21610.428
  # row 113
$ws.Range("H113").Value = 31409.818
$ws.Range("I113").Value = 6312.923
$ws.Range("K113").Value = 6312.923
$ws.Range("M113").Value = -4142.923
  # row 119
$ws.Range("H119").Value = 40000
$ws.Range("J119").Value = 40000
$ws.Range("L119").Value = 40000
$ws.Range("N119").Value = -49676

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
  # row 9
$ws.Range("H9").Value = 517.4
$ws.Range("I9").Value = 517.4
$ws.Range("K9").Value = 517.4
$ws.Range("M9").Value = -293.4
  # row 38
$ws.Range("H38").Value = 31015
$ws.Range("I38").Value = 22030
$ws.Range("J38").Value = 40000
$ws.Range("K38").Value = 22030
$ws.Range("L38").Value = 40000
$ws.Range("M38").Value = -21620
$ws.Range("N38").Value = -40820
  # row 46
$ws.Range("H46").Value = 3506.9524
$ws.Range("J46").Value = 3681.4211
$ws.Range("L46").Value = 3681.4211
$ws.Range("N46").Value = -4057.4211
  # row 61
$ws.Range("H61").Value = 3705.4285
$ws.Range("I61").Value = 1360.3572
$ws.Range("J61").Value = 8395.571
$ws.Range("K61").Value = 1360.3572
$ws.Range("L61").Value = 8395.571
$ws.Range("M61").Value = -1158.3572
$ws.Range("N61").Value = -8799.571
  # row 113
$ws.Range("H113").Value = 3705.4285
$ws.Range("I113").Value = 1360.3572
$ws.Range("J113").Value = 8395.571
$ws.Range("K113").Value = 1360.3572
$ws.Range("L113").Value = 8395.571
$ws.Range("M113").Value = 809.6428000000001
$ws.Range("N113").Value = -12735.571
  # row 132
$ws.Range("H132").Value = 2562
$ws.Range("I132").Value = 3460
$ws.Range("J132").Value = 1215
$ws.Range("K132").Value = 10380
$ws.Range("L132").Value = 3645
$ws.Range("M132").Value = -7850
$ws.Range("N132").Value = -8705
  # row 136
$ws.Range("H136").Value = 9187.4
$ws.Range("I136").Value = 9044.733
$ws.Range("K136").Value = 27134.199
$ws.Range("M136").Value = -24584.199

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
  # row 69
$ws.Range("H69").Value = 10000
$ws.Range("J69").Value = 10000
$ws.Range("L69").Value = 10000
$ws.Range("N69").Value = -11498
  # row 72
$ws.Range("H72").Value = 10000
$ws.Range("J72").Value = 10000
$ws.Range("L72").Value = 30000
$ws.Range("N72").Value = -37488
  # row 81
$ws.Range("H81").Value = 4538.25
$ws.Range("J81").Value = 4986.625
$ws.Range("L81").Value = 9973.25
$ws.Range("N81").Value = -12095.25
  # row 84
$ws.Range("H84").Value = 4538.25
$ws.Range("J84").Value = 4986.625
$ws.Range("L84").Value = 49866.25
$ws.Range("N84").Value = -60474.25
  # row 92
$ws.Range("H92").Value = 29949.75
$ws.Range("J92").Value = 29949.75
$ws.Range("L92").Value = 29949.75
$ws.Range("N92").Value = -34941.75
  # row 132
$ws.Range("H132").Value = 4298.143
$ws.Range("I132").Value = 4298.143
$ws.Range("K132").Value = 12894.429
$ws.Range("M132").Value = -10364.429
  # row 136
$ws.Range("H136").Value = 2476.6667
$ws.Range("I136").Value = 2169.5625
$ws.Range("K136").Value = 6508.6875
$ws.Range("M136").Value = -3958.6875
